$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'63.925.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.89%  "
$ws.Range("D3").Value = "'3.113.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.07%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'609.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").Value = "'145.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.60%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'3.113.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.08%  "
$ws.Range("D9").Value = "'0.520"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.66%  "
$ws.Range("E10").Value = "  -8.17%  "
$ws.Range("E11").Value = "  -9.72%  "
$ws.Range("D12").Value = "'0.470"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.85%  "
$ws.Range("D13").Value = "'0.0000250"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.39%  "
$ws.Range("D14").Value = "'35.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -10.20%  "
$ws.Range("D15").Value = "'3.618.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.21%  "
$ws.Range("E16").Value = "  +1.22%  "
$ws.Range("D17").Value = "'63.924.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.96%  "
$ws.Range("D18").Value = "'3.107.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.30%  "
$ws.Range("D19").Value = "'6.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.21%  "
$ws.Range("D20").Value = "'477.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.78%  "
$ws.Range("D21").Value = "'14.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.92%  "
$ws.Range("D22").Value = "'0.699"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.67%  "
$ws.Range("D23").Value = "'7.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.51%  "
$ws.Range("D24").Value = "'13.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.96%  "
$ws.Range("D25").Value = "'83.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.94%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -9.20%  "
$ws.Range("D28").Value = "'8.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.17%  "
$ws.Range("E29").Value = "  -12.24%  "
$ws.Range("E30").Value = "  -10.66%  "
$ws.Range("D31").Value = "'6.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.87%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("E33").Value = "  -6.59%  "
$ws.Range("D34").Value = "'26.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.05%  "
$ws.Range("D35").Value = "'1.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.10%  "
$ws.Range("D36").Value = "'5.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.91%  "
$ws.Range("D37").Value = "'52.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.82%  "
$ws.Range("D38").Value = "'0.0₃0748"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.79%  "
$ws.Range("D39").Value = "'461.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.06%  "
$ws.Range("D40").Value = "'2.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -14.76%  "
$ws.Range("E41").Value = "  -8.04%  "
$ws.Range("E42").Value = "  -8.77%  "
$ws.Range("D43").Value = "'8.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.55%  "
$ws.Range("D44").Value = "'2.841.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.267"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.00%  "
$ws.Range("E46").Value = "  -13.79%  "
$ws.Range("D47").Value = "'2.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "'26.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.55%  "
$ws.Range("E50").Value = "  -5.40%  "
$ws.Range("D51").Value = "'118.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.49%  "
